$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 191270.432856
$ws.Range("D2").Value = 90.021896

$ws.Range("B3").Value = 80486.17425500001
$ws.Range("D3").Value = 18.940507
$ws.Range("E3").Value = 0

$ws.Range("B4").Value = 705403.760317
$ws.Range("C4").Value = 332

$ws.Range("G5").Value = -20.245188
$ws.Range("H5").Value = -36.957398
$ws.Range("I5").Value = -3.532978
$ws.Range("J5").Value = 0.012793

$ws.Range("G6").Value = 13.89398
$ws.Range("H6").Value = -3.659507
$ws.Range("I6").Value = 31.447467
$ws.Range("J6").Value = 0.151024

$ws.Range("G7").Value = 34.139168
$ws.Range("H7").Value = 20.940118
$ws.Range("I7").Value = 47.338218
$ws.Range("J7").Value = 0
